{"js": "// Cover letter rewrite (\"First version of Skills\"):\n//  - Intro paragraph (\"Hello,\") becomes the opening bio sentence, and the\n//    `_GoBack` bookmark moves from the end of the old bio paragraph to the\n//    very start of this (now merged) paragraph.\n//  - The old bio paragraph's runs collapse into the new \"I built a site...\"\n//    sentence.\n//  - The \"After this, ... four children.\" paragraph is removed entirely.\n//  - The \"I am building a site...\" paragraph becomes \"You will find a\n//    series of questions...\"\n//  - The blank paragraph + \"I share the link...\" paragraph + blank paragraph\n//    that used to follow are removed (the link / remaining paragraphs stay).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph indices in the original document:\n// 0 \"Hello,\"\n// 1 \"  My name is Juan Pablo Nicotra and i am 41 years old. ... be able to apply.\" (+ _GoBack bookmark)\n// 2 \"After this, my intention ... four children.\"\n// 3 \"I am building a site with my resume but it is not finished yet.\"\n// 4 \"\" (empty)\n// 5 \"I share the link with you, ... about my profile.\"\n// 6 \"\" (empty)\n// 7 \"https://jpnicotra.github.io/myCv/\"\n// ... (unchanged tail)\n\nconst pHello = paragraphs.items[0];\nconst pBio = paragraphs.items[1];\nconst pAfterThis = paragraphs.items[2];\nconst pBuildingSite = paragraphs.items[3];\nconst pEmptyAfterSite = paragraphs.items[4];\nconst pShareLink = paragraphs.items[5];\nconst pEmptyAfterShare = paragraphs.items[6];\n\n// Remove the old `_GoBack` bookmark wherever it currently lives (end of the\n// bio paragraph) before re-creating it at the new location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// \"Hello,\" -> new opening sentence, with the bookmark re-inserted at the\n// very start of the paragraph (before the run).\npHello.insertText(\n  \"My name is Juan Pablo Nicotra. I am 41 years old. I am from Argentina and in the process of obtaining Italian citizenship. I need to finish some documents and travel to Italy to complete the process.\",\n  \"Replace\"\n);\npHello.getRange(\"Start\").insertBookmark(\"_GoBack\");\n\n// Old bio paragraph's text becomes the \"I built a site...\" sentence.\npBio.insertText(\n  \"I built a site with my resume and some other aspects of my life that probably should need in this selection process.\",\n  \"Replace\"\n);\n\n// The \"After this, ... four children.\" paragraph is dropped completely.\npAfterThis.delete();\n\n// \"I am building a site...\" -> \"You will find a series of questions...\"\npBuildingSite.insertText(\n  \"You will find a series of questions and answers to learn more about my profile.\",\n  \"Replace\"\n);\n\n// Drop the blank paragraph, the \"I share the link...\" paragraph, and the\n// blank paragraph that followed it.\npEmptyAfterSite.delete();\npShareLink.delete();\npEmptyAfterShare.delete();\n\nawait context.sync();\n", "ps1": "# Cover letter rewrite (\"First version of Skills\"):\n#  - Intro paragraph (\"Hello,\") becomes the opening bio sentence, and the\n#    `_GoBack` bookmark moves from the end of the old bio paragraph to the\n#    very start of this (now merged) paragraph.\n#  - The old bio paragraph's runs collapse into the new \"I built a site...\"\n#    sentence.\n#  - The \"After this, ... four children.\" paragraph is removed entirely.\n#  - The \"I am building a site...\" paragraph becomes \"You will find a\n#    series of questions...\"\n#  - The blank paragraph + \"I share the link...\" paragraph + blank paragraph\n#    that used to follow are removed (the link / remaining paragraphs stay).\n\n$d = $word.ActiveDocument\n\n# Original (1-based) paragraph numbering:\n#  1 \"Hello,\"\n#  2 \"  My name is Juan Pablo Nicotra and i am 41 years old. ... be able to apply.\" (+ _GoBack bookmark)\n#  3 \"After this, my intention ... four children.\"\n#  4 \"I am building a site with my resume but it is not finished yet.\"\n#  5 \"\" (empty)\n#  6 \"I share the link with you, ... about my profile.\"\n#  7 \"\" (empty)\n#  8 \"https://jpnicotra.github.io/myCv/\"\n#  ... (unchanged tail)\n\n# Move the `_GoBack` bookmark: delete it from wherever it currently sits\n# (end of paragraph 2) before re-inserting it at the start of paragraph 1.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Paragraph 1 (\"Hello,\") -> new opening sentence.\n$p1 = $d.Paragraphs.Item(1)\n$r1 = $p1.Range\n$r1.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n$r1.Text = \"My name is Juan Pablo Nicotra. I am 41 years old. I am from Argentina and in the process of obtaining Italian citizenship. I need to finish some documents and travel to Italy to complete the process.\"\n\n# Re-insert `_GoBack` at the very start of paragraph 1.\n$startRange = $d.Paragraphs.Item(1).Range.Duplicate()\n$startRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $startRange) | Out-Null\n\n# Paragraph 2 (old bio paragraph) -> \"I built a site...\" sentence.\n$p2 = $d.Paragraphs.Item(2)\n$r2 = $p2.Range\n$r2.MoveEnd(1, -1) | Out-Null\n$r2.Text = \"I built a site with my resume and some other aspects of my life that probably should need in this selection process.\"\n\n# Paragraph 3 (\"After this, ... four children.\") is dropped completely.\n$d.Paragraphs.Item(3).Range.Delete()\n\n# Paragraph 4 (\"I am building a site...\") -> \"You will find a series of\n# questions...\" (still at index 3 after the previous delete).\n$p4 = $d.Paragraphs.Item(3)\n$r4 = $p4.Range\n$r4.MoveEnd(1, -1) | Out-Null\n$r4.Text = \"You will find a series of questions and answers to learn more about my profile.\"\n\n# Drop the blank paragraph, the \"I share the link...\" paragraph, and the\n# blank paragraph that followed it (all now sitting right after index 3).\n$d.Paragraphs.Item(4).Range.Delete()\n$d.Paragraphs.Item(4).Range.Delete()\n$d.Paragraphs.Item(4).Range.Delete()\n"}
